$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.647.86"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.113.51"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'243.53"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("D6").Value = "'624.53"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'1.16"
$ws.Range("E7").Value = "  +12.21%  "
$ws.Range("D8").Value = "'0.371"
$ws.Range("E8").Value = "  +5.99%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.111.23"
$ws.Range("E10").Value = "  -8.00%  "
$ws.Range("D11").Value = "'0.751"
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").Value = "'35.34"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "90.523.60"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "3.685.41"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "3.076.26"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "'3.82"
$ws.Range("E19").Value = "  +4.45%  "
$ws.Range("D20").Value = "'14.39"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  +7.91%  "
$ws.Range("D23").Value = "'448.21"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'9.12"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").Value = "'5.88"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").Value = "'93.46"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").Value = "'11.92"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.177"
$ws.Range("E30").Value = "  +11.99%  "
$ws.Range("D31").Value = "'0.224"
$ws.Range("E31").Value = "  +12.94%  "
$ws.Range("D32").Value = "'9.13"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +35.55%  "
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +30.20%  "
$ws.Range("E35").Value = "  +5.08%  "
$ws.Range("D36").Value = "'26.61"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").Value = "'7.68"
$ws.Range("E37").Value = "  +10.12%  "
$ws.Range("D38").Value = "'4.21"
$ws.Range("E38").Value = "  +28.94%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'494.43"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "'3.62"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "'1.29"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'0.417"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'22.11"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'157.42"
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "'0.686"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'4.56"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").Value = "'45.01"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  +0.55%  "
